$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '27.661.41'
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.847.23'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -0.79%  '
$ws.Range("E4").Value = '  -0.26%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '314.41'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.62%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.24%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4266'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -3.23%  '
$ws.Range("E8").Value = '  -1.33%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '44.77'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.94%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.07241'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -4.02%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.9008'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -3.89%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '20.70'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -2.91%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '1.822.14'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -3.86%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '6.578'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -1.98%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '5.355'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -2.05%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.06835'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.57%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '1.004'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.17%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '77.68'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -5.61%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.000008832'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -3.36%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.36%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '15.47'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -3.46%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '27.656.01'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -2.46%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '4.962'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -3.14%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '10.67'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.66%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.082.16'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -1.77%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '2.048'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.32%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '154.04'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.47%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '18.18'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.18%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '5.243'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -2.54%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.829'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +5.21%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '111.03'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -3.27%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.08896'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -2.33%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.7716'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -4.06%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '4.562'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -6.69%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '2.918'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.70%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.083'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -7.86%  '
$ws.Range("E37").Value = '  -0.36%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.05379'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -1.66%  '
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '1.096'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -3.09%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '2.956'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -1.73%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.01923'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -1.59%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.5065'
$c.Style = "Normal"
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.1641'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -2.29%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '6.762'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -5.75%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '8.265'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -6.93%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.06642'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -2.07%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '10.39'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -1.89%  '
$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.4715'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -3.87%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '105.18'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -2.50%  '
$ws.Range("E50").Value = '  -0.27%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.639'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -2.85%  '
